$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(6, 6, 7, 4, 8, 3, 9, 3, 7, 5, 6)
$jValues = @(6, 7, 8, 6, 9, 3, 9, 3, 7, 5, 6)

for ($r = 0; $r -lt 11; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
